# CA3/data/draft1.xlsx edit: add a "spectral line color" column to Sheet1,
# append calculated peak values + a small results table, duplicate the
# original Sheet1 layout into a new Sheet3, and touch up selections.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the ORIGINAL Sheet1 (before any of today's edits) into a
#    new Sheet3, with column B header changed to "实验值1/nm".
# ---------------------------------------------------------------------
$s1 = $wb.Worksheets.Item("Sheet1")
$s1.Range("A1:C6").Copy()

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$s3 = $wb.Worksheets.Add($null, $last)
$s3.Name = "Sheet3"
$s3.Range("B1").PasteSpecial(-4104)
$s3.Range("B1").Value = "实验值1/nm"

$s3.Range("D8").Select()

# ---------------------------------------------------------------------
# 2) Sheet1: insert a new first column holding the spectral line colour
#    for each row (shifts old A/B/C -> B/C/D, and fixes up the formulas
#    automatically).
# ---------------------------------------------------------------------
$s1.Columns.Item(1).Insert()

$s1.Range("A1").Value = "谱线颜色"

$s1.Range("A2").Value = "黄色"
$s1.Range("A3").Value = "黄色"
$s1.Range("A4").Value = "绿色"
$s1.Range("A5").Value = "蓝色"
$s1.Range("A6").Value = "紫色"

# style the colour-name cells: bigger SimSun font, centered + wrapped
$colorCells = $s1.Range("A2:A6")
$colorCells.Font.Name = "SimSun"
$colorCells.Font.Size = 12
$colorCells.Font.Color = 0
$colorCells.HorizontalAlignment = -4108
$colorCells.VerticalAlignment = -4108
$colorCells.WrapText = $true

$s1.Rows.Item(2).RowHeight = 15.6
$s1.Rows.Item(3).RowHeight = 15.6
$s1.Rows.Item(4).RowHeight = 15.6
$s1.Rows.Item(5).RowHeight = 15.6
$s1.Rows.Item(6).RowHeight = 15.6

# extra measured/derived values off to the side
$s1.Range("L4").Value = 1000000000

$s1.Range("H5").Value = 436.30720922031827
$s1.Range("I5").Value = 546.2969403898162
$s1.Range("J5").Value = 578.80414613417565

$s1.Range("H6").Value = 0.00000043630720922031825
$s1.Range("I6").Value = 0.00000054629694038981619
$s1.Range("J6").Value = 0.00000057880414613417571

# small "实验值" results table underneath
$s1.Range("B11").Value = "实验值"

$s1.Range("A12").Value = "蓝色"
$s1.Range("B12").Value = 436.30720922031827

$s1.Range("A13").Value = "绿色"
$s1.Range("B13").Value = 546.2969403898162

$s1.Range("A14").Value = "黄色"
$s1.Range("B14").Value = 578.80414613417565

$colorCells2 = $s1.Range("A12:A14")
$colorCells2.Font.Name = "SimSun"
$colorCells2.Font.Size = 12
$colorCells2.Font.Color = 0
$colorCells2.HorizontalAlignment = -4108
$colorCells2.VerticalAlignment = -4108
$colorCells2.WrapText = $true

$s1.Rows.Item(12).RowHeight = 15.6
$s1.Rows.Item(13).RowHeight = 15.6
$s1.Rows.Item(14).RowHeight = 15.6

# ---------------------------------------------------------------------
# 3) Sheet2: tweak the view (top-left + selection) a little.
# ---------------------------------------------------------------------
$s2 = $wb.Worksheets.Item("Sheet2")
$s2.Activate()
$s2.Range("C21").Select()

# ---------------------------------------------------------------------
# 4) Sheet1 stays the active tab, with C12 selected last.
# ---------------------------------------------------------------------
$s1.Activate()
$s1.Range("C12").Select()

Write-Output "edit applied"
